# Update crypto price/volume data per latest scrape (GitHub Actions bot)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.257.86'
$ws.Range("E2").Value = '  +0.50%  '

$ws.Range("D3").Value = '1.858.14'
$ws.Range("E3").Value = '  +0.39%  '

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = "'0.7051"
$ws.Range("E5").Value = '  +1.82%  '

$ws.Range("D6").Value = "'238.34"
$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = "'0.08015"
$ws.Range("E8").Value = '  +3.64%  '

$ws.Range("E9").Value = '  -0.34%  '

$ws.Range("D10").Value = "'23.51"
$ws.Range("E10").Value = '  +1.12%  '

$ws.Range("D11").Value = "'0.08187"
$ws.Range("E11").Value = '  +0.85%  '

$ws.Range("D12").Value = '1.910.45'
$ws.Range("E12").Value = '  +2.92%  '

$ws.Range("D13").Value = "'5.197"
$ws.Range("E13").Value = '  -0.09%  '

$ws.Range("D14").Value = "'0.7066"
$ws.Range("E14").Value = '  -2.55%  '

$ws.Range("D15").Value = "'89.69"
$ws.Range("E15").Value = '  +0.76%  '

$ws.Range("D16").Value = '29.249.24'
$ws.Range("E16").Value = '  +0.40%  '

$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").Value = "'5.833"
$ws.Range("E17").Value = '  +1.54%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = "'0.000007925"
$ws.Range("E18").Value = '  +1.07%  '

$ws.Range("D19").Value = "'13.28"
$ws.Range("E19").Value = '  +0.73%  '

$ws.Range("D20").Value = "'238.32"
$ws.Range("E20").Value = '  +1.00%  '

$ws.Range("D21").Value = "'0.9988"

$ws.Range("D22").Value = '2.102.66'
$ws.Range("E22").Value = '  +0.07%  '

$ws.Range("D23").Value = "'1.000"

$ws.Range("D24").Value = "'7.475"
$ws.Range("E24").Value = '  -1.67%  '

$ws.Range("D25").Value = "'162.86"
$ws.Range("E25").Value = '  +0.96%  '

$ws.Range("D26").Value = "'8.874"
$ws.Range("E26").Value = '  -1.08%  '

$ws.Range("E27").Value = '  +0.74%  '

$ws.Range("E28").Value = '  +0.29%  '

$ws.Range("D29").Value = "'1.931"
$ws.Range("E29").Value = '  -1.91%  '

$ws.Range("D30").Value = "'1.430"
$ws.Range("E30").Value = '  +2.20%  '

$ws.Range("D31").Value = "'1.475"
$ws.Range("E31").Value = '  -0.81%  '

$ws.Range("D32").Value = "'4.371"
$ws.Range("E32").Value = '  -2.95%  '

$ws.Range("D33").Value = "'4.027"
$ws.Range("E33").Value = '  +0.73%  '

$ws.Range("D34").Value = "'0.05202"
$ws.Range("E34").Value = '  -0.13%  '

$ws.Range("D35").Value = "'1.161"
$ws.Range("E35").Value = '  -1.60%  '

$ws.Range("D36").Value = "'0.7169"
$ws.Range("E36").Value = '  +1.78%  '

$ws.Range("D37").Value = "'1.001"
$ws.Range("E37").Value = '  -2.48%  '

$ws.Range("D38").Value = "'2.668"

$ws.Range("D39").Value = "'0.01859"
$ws.Range("E39").Value = '  +0.26%  '

$ws.Range("D40").Value = "'2.730"
$ws.Range("E40").Value = '  +2.20%  '

$ws.Range("D41").Value = "'0.9373"
$ws.Range("E41").Value = '  +2.70%  '

$ws.Range("D42").Value = '1.137.93'
$ws.Range("E42").Value = '  +3.89%  '

$ws.Range("D43").Value = "'5.985"
$ws.Range("E43").Value = '  -0.54%  '

$ws.Range("D44").Value = "'0.4268"
$ws.Range("E44").Value = '  -0.07%  '

$ws.Range("D45").Value = "'70.66"
$ws.Range("E45").Value = '  +0.07%  '

$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").Value = "'102.88"
$ws.Range("E47").Value = '  -0.21%  '

$ws.Range("D48").Value = "'0.5288"
$ws.Range("E48").Value = '  -4.56%  '

$ws.Range("E49").Value = '  -0.04%  '

$ws.Range("D50").Value = '2.003.96'
$ws.Range("E50").Value = '  +0.15%  '

$ws.Range("D51").Value = "'9.174"
$ws.Range("E51").Value = '  +0.17%  '
